$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 207 ---
# A207 already holds the shared string "25-10-2021" and is left untouched.
$ws.Range("B207").Value = 68.96
$ws.Range("C207").Value = 15.24

# --- Append new rows 208-211 (unambiguous day-of-month, safe as plain text) ---
$ws.Range("A208").Value = "26-10-2021"
$ws.Range("B208").Value = 66.7
$ws.Range("C208").Value = 15.98

$ws.Range("A209").Value = "27-10-2021"
$ws.Range("B209").Value = 70.06
$ws.Range("C209").Value = 16.98

$ws.Range("A210").Value = "28-10-2021"
$ws.Range("B210").Value = 71.96
$ws.Range("C210").Value = 16.53

$ws.Range("A211").Value = "29-10-2021"
$ws.Range("B211").Value = 75.45
$ws.Range("C211").Value = 16.26

# --- Rows 212-213 use dates whose day-of-month is <= 12, which Excel's   ---
# --- smart typing would otherwise auto-convert into a real date value.   ---
# --- Build the literal text with a formula (never re-parsed) in a spare ---
# --- helper cell, copy it, then paste-special just the value so the     ---
# --- destination keeps plain text with the default (unstyled) format.   ---
$helper = $ws.Range("ZZ1")

$helper.Formula = "=""01-11-2021"""
$helper.Copy()
$ws.Range("A212").PasteSpecial(-4163)
$helper.ClearContents()

$ws.Range("B212").Value = 78.34
$ws.Range("C212").Value = 16.41

$helper.Formula = "=""02-11-2021"""
$helper.Copy()
$ws.Range("A213").PasteSpecial(-4163)
$helper.ClearContents()

$ws.Range("C213").Value = 16.42
